# New weekly record: insert a new row at row 4 (shifting the existing
# rows 4-39 down to 5-40) and populate it with the latest "Rabanito"
# price observation (week of 2022-11-17).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("4:4").Insert()

$ws.Cells.Item(4, 1).Value = 6
$ws.Cells.Item(4, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(4, 3).Value = "Metropolitana"
$ws.Cells.Item(4, 4).Value = 44882
$ws.Cells.Item(4, 5).Value = 13
$ws.Cells.Item(4, 6).Value = 300000001
$ws.Cells.Item(4, 7).Value = "Rabanito"
$ws.Cells.Item(4, 8).Value = "Sin especificar"
$ws.Cells.Item(4, 9).Value = "Primera"
$ws.Cells.Item(4, 10).Value = 7900
$ws.Cells.Item(4, 11).Value = 3000
$ws.Cells.Item(4, 12).Value = 3000
$ws.Cells.Item(4, 13).Value = 3000
$ws.Cells.Item(4, 14).Value = "`$/cien unidades (volumen en unidades)"
$ws.Cells.Item(4, 15).Value = "Provincia de Chacabuco"
$ws.Cells.Item(4, 16).Value = 30
$ws.Cells.Item(4, 17).Value = 100
$ws.Cells.Item(4, 18).Value = "Hortaliza"
